# Word COM-interop script: remove the "logical()" factor significance
# columns from the "Profilaxia (%)" row and replace them with a new
# breakdown of rows (0..4) showing n (%) per sex, inserted right after
# the "Profilaxia (%)" row and before the "Dabigatrana (%)" row.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CompactCenteredText($cell, [string]$text, [bool]$bold) {
    $cell.Range.Text = $text
    $cell.Range.Paragraphs.Item(1).Style = "Compact"
    $cell.Range.ParagraphFormat.Alignment = 1
    if ($text.Length -gt 0) {
        $r = $d.Range($cell.Range.Start, $cell.Range.End - 1)
        $r.Bold = [int]$bold
    }
}

function Clear-CompactCell($cell) {
    $r = $d.Range($cell.Range.Start, $cell.Range.End - 1)
    $r.Delete()
    $cell.Range.ParagraphFormat.Alignment = 0
    $cell.Range.Paragraphs.Item(1).Style = "Compact"
}

# ---------------------------------------------------------------------
# 1) "Profilaxia (%)" row (row 3): the sex-specific counts and the old
#    p-value are gone; only the new (smaller) p-value 0.010 remains.
# ---------------------------------------------------------------------
$profRow = $t.Rows.Item(3)
Clear-CompactCell $profRow.Cells.Item(2)
Clear-CompactCell $profRow.Cells.Item(3)
Set-CompactCenteredText $profRow.Cells.Item(4) "0.010" $false

# ---------------------------------------------------------------------
# 2) Insert five new rows (labelled 0..4) right before "Dabigatrana (%)"
#    (currently row 4), each row needs to end up with 5 cells: label,
#    two counts, and two blank "Compact" cells.
# ---------------------------------------------------------------------
$newRowsData = @(
    @("4", "2 ( 0.9)",   "0 ( 0.0)"),
    @("3", "16 ( 7.0)",  "15 (10.9)"),
    @("2", "154 (67.8)", "71 (51.4)"),
    @("1", "50 (22.0)",  "44 (31.9)"),
    @("0", "5 ( 2.2)",   "8 ( 5.8)")
)

$refRow = $t.Rows.Item(4)
foreach ($rowData in $newRowsData) {
    $newRow = $t.Rows.Add($refRow)

    $c1 = $newRow.Cells.Item(1)
    $c1.Split(1, 5)
    $c2 = $newRow.Cells.Item(2)
    $c2.Split(1, 5)
    $c3 = $newRow.Cells.Item(3)
    $c3.Split(1, 5)
    $c4 = $newRow.Cells.Item(4)
    $c4.Split(1, 5)

    Set-CompactCenteredText $newRow.Cells.Item(1) $rowData[0] $true
    Set-CompactCenteredText $newRow.Cells.Item(2) $rowData[1] $false
    Set-CompactCenteredText $newRow.Cells.Item(3) $rowData[2] $false

    $c4b = $newRow.Cells.Item(4)
    $c4b.Range.ParagraphFormat.Alignment = 0
    $c4b.Range.Paragraphs.Item(1).Style = "Compact"

    $c5b = $newRow.Cells.Item(5)
    $c5b.Range.ParagraphFormat.Alignment = 0
    $c5b.Range.Paragraphs.Item(1).Style = "Compact"
}

Write-Host "Done"
